$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D:E").Insert()

# Copy number-format/style from the (now-shifted) old D:E columns - now at F:G -
# into the newly inserted D:E columns, in the three contiguous blocks that
# actually contained data in the original sheet (rows 7-35, 38-77, 80-102).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new quarters of data (Dec-2018 / Sep-2018) in columns D and E
$ws.Range("D7").Value = 43465; $ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 14500; $ws.Range("E8").Value = 7300
$ws.Range("D9").Value = 400; $ws.Range("E9").Value = 200
$ws.Range("D10").Value = 14100; $ws.Range("E10").Value = 7100
$ws.Range("D12").Value = 26500; $ws.Range("E12").Value = 28500
$ws.Range("D13").Value = 0; $ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0; $ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0; $ws.Range("E15").Value = 0
$ws.Range("D17").Value = 66400; $ws.Range("E17").Value = 68600
$ws.Range("D18").Value = -51900; $ws.Range("E18").Value = -61300
$ws.Range("D20").Value = 1200; $ws.Range("E20").Value = -23300
$ws.Range("D21").Value = -49900; $ws.Range("E21").Value = -84000
$ws.Range("D22").Value = 800; $ws.Range("E22").Value = 800
$ws.Range("D23").Value = -51500; $ws.Range("E23").Value = -85400
$ws.Range("D24").Value = 0; $ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0; $ws.Range("E25").Value = 0
$ws.Range("D26").Value = -51500; $ws.Range("E26").Value = -85400
$ws.Range("D27").Value = -51500; $ws.Range("E27").Value = -85400
$ws.Range("D28").Value = 0; $ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0; $ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0; $ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0; $ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1200; $ws.Range("E32").Value = 23300
$ws.Range("D33").Value = -51500; $ws.Range("E33").Value = -85400
$ws.Range("D34").Value = 0; $ws.Range("E34").Value = 0
$ws.Range("D35").Value = -51500; $ws.Range("E35").Value = -85400
$ws.Range("D38").Value = 43465; $ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 202800; $ws.Range("E41").Value = 235000
$ws.Range("D42").Value = 0; $ws.Range("E42").Value = 1000
$ws.Range("D43").Value = 2700; $ws.Range("E43").Value = 2000
$ws.Range("D44").Value = 10100; $ws.Range("E44").Value = 5600
$ws.Range("D45").Value = 4500; $ws.Range("E45").Value = 3300
$ws.Range("D46").Value = 220200; $ws.Range("E46").Value = 246800
$ws.Range("D47").Value = "NA"; $ws.Range("E47").Value = "NA"
$ws.Range("D48").Value = 60500; $ws.Range("E48").Value = 58400
$ws.Range("D49").Value = 0; $ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0; $ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0; $ws.Range("E51").Value = 0
$ws.Range("D52").Value = 4300; $ws.Range("E52").Value = 4000
$ws.Range("D53").Value = 0; $ws.Range("E53").Value = 0
$ws.Range("D54").Value = 285000; $ws.Range("E54").Value = 309200
$ws.Range("D57").Value = 12400; $ws.Range("E57").Value = 7100
$ws.Range("D58").Value = 0; $ws.Range("E58").Value = 0
$ws.Range("D59").Value = 38400; $ws.Range("E59").Value = 27900
$ws.Range("D60").Value = 50800; $ws.Range("E60").Value = 35000
$ws.Range("D61").Value = 0; $ws.Range("E61").Value = 0
$ws.Range("D62").Value = 6500; $ws.Range("E62").Value = 5600
$ws.Range("D63").Value = 0; $ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0; $ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0; $ws.Range("E65").Value = 0
$ws.Range("D66").Value = 57200; $ws.Range("E66").Value = 40600
$ws.Range("D68").Value = 0; $ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0; $ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0; $ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0; $ws.Range("E71").Value = 0
$ws.Range("D72").Value = -696400; $ws.Range("E72").Value = -645000
$ws.Range("D73").Value = 0; $ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0; $ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0; $ws.Range("E75").Value = 0
$ws.Range("D76").Value = 227800; $ws.Range("E76").Value = 268600
$ws.Range("D77").Value = 0; $ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465; $ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -51500; $ws.Range("E81").Value = -85400
$ws.Range("D83").Value = 700; $ws.Range("E83").Value = 600
$ws.Range("D84").Value = 0; $ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0; $ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0; $ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0; $ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0; $ws.Range("E88").Value = 0
$ws.Range("D89").Value = -31500; $ws.Range("E89").Value = -43300
$ws.Range("D91").Value = -1900; $ws.Range("E91").Value = -6400
$ws.Range("D92").Value = 0; $ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0; $ws.Range("E93").Value = 0
$ws.Range("D94").Value = -900; $ws.Range("E94").Value = 8100
$ws.Range("D96").Value = 0; $ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0; $ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0; $ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0; $ws.Range("E99").Value = 0
$ws.Range("D100").Value = 300; $ws.Range("E100").Value = -500
$ws.Range("D101").Value = 0; $ws.Range("E101").Value = 0
$ws.Range("D102").Value = -32100; $ws.Range("E102").Value = -35700

